$wb = $excel.ActiveWorkbook

$wsTriangle = $wb.Worksheets.Item("triangle")
$wsLdfs = $wb.Worksheets.Item("LDFs")

# --- triangle sheet: row 3 height 43.5 -> 72.5 ---
$wsTriangle.Rows.Item(3).RowHeight = 72.5

# --- LDFs sheet: B28 goes from 1 to 11, and the "+1" countdown chain becomes a "-1" countdown chain ---
$wsLdfs.Range("B28").Value = 11
$wsLdfs.Range("C28").Formula = "=+B28-1"
$wsLdfs.Range("D28:L28").FormulaR1C1 = "=+RC[-1]-1"

# --- LDFs sheet: reverse the order of the static LDF table B36:B46 to match B15:L15 reversed ---
$srcValues = $wsLdfs.Range("B15:L15").Value2
$n = 11
for ($i = 0; $i -lt $n; $i++) {
    $val = $srcValues[1, $n - $i]
    $destRow = 36 + $i
    $wsLdfs.Cells.Item($destRow, 2).Value = $val
}

# --- LDFs sheet: selection moves from A1:M26 to single active cell B33 ---
$wsLdfs.Activate()
$wsLdfs.Range("B33").Select()
